$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing worker's "Periodo Mora" (2506 -> 2507)
$ws.Range("E16").Value = "2507"

# Insert a new row below the existing worker row (16), copying its formatting,
# to host the new worker's data. This also shifts the two signature-line rows
# (formerly 21 and 22) down by one, to 22 and 23.
$ws.Rows.Item(16).Copy()
$ws.Rows.Item(17).Insert()

# Fill in the new worker's data on the newly inserted row 17
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047511606"
$ws.Range("D17").Value = "CAMILA MARCELA PRIMERA GUERRERO"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Update the summary totals: total "Valor Mora" and worker count now reflect
# both workers.
$ws.Range("E11").Value = 113880
$ws.Range("C13").Value = 2

# Widen column D so the longer worker name fits (mirrors Excel's bestFit
# recalculation after the longer name was entered).
$ws.Columns.Item(4).ColumnWidth = 35.83

Write-Output "done"
